$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update the existing task list (rows 2-5): new task text ----
$ws.Range("A2").Value = "Implement feature 1 (Ricardo)"
$ws.Range("A3").Value = "Implement feature 2 (James, João, Francisco, Iago)"
$ws.Range("A4").Value = "Upload sprints to github (Todo mundo)"
$ws.Range("A5").Value = "Take metrics (João)"

# Remove the old 5th task row entirely
$ws.Range("A6").ClearContents()

# Shade the updated rows (A2:D5) with the new light-blue fill
$ws.Range("A2:D5").Interior.Color = 15917529

# ---- Build a second Kanban header block at row 8, reusing row 1's formats ----
$ws.Range("A1:D1").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Todo"
$ws.Range("B8").Value = "In Progress"
$ws.Range("C8").Value = "Reviewing"
$ws.Range("D8").Value = "Done"

# ---- Second block body (rows 9-12): same fill, values now live in column D ----
$ws.Range("A9:D12").Interior.Color = 15917529
$ws.Range("D9").Value = "Implement feature 1 (Ricardo)"
$ws.Range("D10").Value = "Implement feature 2 (James, João, Francisco, Iago)"
$ws.Range("D11").Value = "Upload sprints to github (Todo mundo)"
$ws.Range("D12").Value = "Take metrics (João)"

# ---- Column widths ----
$ws.Columns.Item(1).ColumnWidth = 56.66666666666667
$ws.Columns.Item(4).ColumnWidth = 56.5

# ---- Selection ----
$ws.Range("B20").Select()
